$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update id schema effective date (effective_from) and cr_dtimes to the
# new effective date/time, and switch the date/time number format used
# by these cells from a time-only format to a full date + time format.
$newDateSerial = 45079.634583333333

# effective_from
$ws.Range("H2").Value = $newDateSerial
$ws.Range("H2").NumberFormat = "m/d/yy h:mm"

# cr_dtimes
$ws.Range("L2").Value = $newDateSerial
$ws.Range("L2").NumberFormat = "m/d/yy h:mm"

$ws.Range("L4").Select() | Out-Null
